$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.700.60'

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.895.82'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.30%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -1.08%  '

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.36%  '

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.99%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4865'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.73%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3793'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.60%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07323'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.71%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -2.78%  '

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -2.25%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07660'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.75%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.874.61'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.21%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.480'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.41%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.596'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.22%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.16'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.13%  '

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.08%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008777'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.96%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.92%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.580.20'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.99%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.49'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -2.48%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.120'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.131.88'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.47%  '

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.99%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.909'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.00%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '153.69'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.40%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.37'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.04%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.150'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +4.73%  '

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.43%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.866'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.28%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08890'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.05%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.203'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -4.00%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.221'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.77%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7656'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.60%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.629'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.75%  '

$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02036'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.38%  '

$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.524'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -7.72%  '

$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.095'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.07%  '

$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05271'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.02%  '

$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5470'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -2.47%  '

$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.974'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.85%  '

$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.888'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.30%  '

$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.499'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.67%  '

$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1518'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.64%  '

$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '111.71'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +5.96%  '

$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.63'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.88%  '

$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4782'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.02%  '

$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.09%  '

$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.632'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.90%  '

$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '67.14'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.32%  '

$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06050'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.24%  '
